$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Insert a new row above row 26, shifting existing rows 26-74 down to 27-75
$ws.Rows.Item(26).Insert()

# Fill in the new row's data (name, label, valueType)
$ws.Cells.Item(26, 2).Value = "J38c"
$ws.Cells.Item(26, 3).Value = "SACANA FGS J38c ""Sweetened cheese and curd"""
$ws.Cells.Item(26, 4).Value = "decimal"

# Update the selection/active cell to D26 and clear any frozen top-left cell
$ws.Range("D26").Select()
